# --- Rename the worksheet tab (matches workbook.xml <sheet name="..."> change) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "repayment_20250901_20250920 (2)"

# --- Updated collector data rows (row -> D,E,F,G,H,I,J,K,L) ---
# E, F, G, K, L are stored as literal text strings in the workbook
# (e.g. "41,259,415.00"), not as numbers, so they must be written as text.
$rows = @(
    @{ R=2;  D=58; E="41,259,415.00"; F="389,618,385.00"; G="10.59"; H=22.542999999999999; I=327; J=20; K="9.21"; L="6.12" },
    @{ R=3;  D=52; E="36,671,431.00"; F="383,088,000.00"; G="9.57";  H=26.824999999999999; I=330; J=13; K="2.49"; L="3.94" },
    @{ R=4;  D=46; E="26,880,737.00"; F="395,616,117.00"; G="6.79";  H=17.329000000000001; I=331; J=19; K="5.03"; L="5.74" },
    @{ R=5;  D=61; E="33,400,674.00"; F="358,267,731.00"; G="9.32";  H=23.777000000000001; I=297; J=17; K="5.06"; L="5.72" },
    @{ R=6;  D=48; E="32,500,458.00"; F="370,220,959.00"; G="8.78";  H=17.896000000000001; I=329; J=10; K="2.79"; L="3.04" },
    @{ R=7;  D=71; E="41,871,931.00"; F="406,214,480.00"; G="10.31"; H=13.56;               I=329; J=15; K="3.98"; L="4.56" },
    @{ R=8;  D=59; E="45,786,096.00"; F="381,376,141.00"; G="12.01"; H=27.533000000000001; I=324; J=17; K="3.53"; L="5.25" },
    @{ R=9;  D=44; E="38,639,329.00"; F="427,498,587.00"; G="9.04";  H=15.49;               I=332; J=9;  K="4.32"; L="2.71" },
    @{ R=10; D=30; E="18,923,984.00"; F="341,191,913.00"; G="5.55";  H=16.376000000000001; I=264; J=7;  K="2.85"; L="2.65" },
    @{ R=11; D=52; E="48,209,748.00"; F="383,910,226.00"; G="12.56"; H=15.321999999999999; I=327; J=22; K="9.46"; L="6.73" },
    @{ R=12; D=49; E="30,381,661.00"; F="376,834,676.00"; G="8.06";  H=23.04;               I=331; J=17; K="4.74"; L="5.14" },
    @{ R=13; D=42; E="32,908,908.00"; F="390,335,562.00"; G="8.43";  H=19.608000000000001; I=325; J=11; K="3.55"; L="3.38" },
    @{ R=14; D=56; E="48,655,390.00"; F="376,779,839.00"; G="12.91"; H=12.061999999999999; I=326; J=13; K="9.89"; L="3.99" },
    @{ R=15; D=43; E="29,368,835.00"; F="390,760,194.00"; G="7.52";  H=15.43;               I=329; J=12; K="3.32"; L="3.65" },
    @{ R=16; D=51; E="34,919,292.00"; F="383,189,468.00"; G="9.11";  H=12.459;              I=331; J=12; K="2.64"; L="3.63" },
    @{ R=17; D=49; E="35,780,349.00"; F="369,862,031.00"; G="9.67";  H=24.027000000000001; I=326; J=14; K="3.66"; L="4.29" },
    @{ R=18; D=33; E="26,408,642.00"; F="320,145,256.00"; G="8.25";  H=14.779;              I=230; J=4;  K="2.47"; L="1.74" }
)

# Force the text-valued columns (E,F,G,K,L) to be written as literal strings
# instead of being auto-converted to numbers by the COM value setter: set the
# range to Text format first, assign the values, then clear the formatting
# again so the cells end up with no explicit style (matching a fresh default
# cell) while keeping their content as shared-string text.
$textRange = $ws.Range("E2:G18")
$textRange.NumberFormat = "@"
$textRange2 = $ws.Range("K2:L18")
$textRange2.NumberFormat = "@"

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
}

$textRange.ClearFormats()
$textRange2.ClearFormats()
